$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46: politeness_score (B46) was stored as text "3" - normalize it to the
# numeric value 3, matching the rest of the column.
$ws.Range("B46").Value = 3

# Row 47: brand new annotation row appended after row 46.
$ws.Range("A47").Value = "Sunsi Wu"

# B47 keeps its "4" as text (inline string), not a number, so force text
# formatting before assignment and then drop back to the Normal style so no
# stray number-format/style survives on the cell.
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "4"
$ws.Range("B47").Style = "Normal"

$ws.Range("C47").Value = "will"
$ws.Range("D47").Value = "FBK"
$ws.Range("E47").Value = "WRI"
$ws.Range("F47").Value = "92b80f86-ee70-4a78-8469-1a9c33b052ed"
$ws.Range("G47").Value = "7Y52YHDS2X7ae_annotated.xlsx"
$ws.Range("H47").Value = "We will include this description in the new version of the paper."
